$wb = $excel.ActiveWorkbook

# Update the "Ready for handoff" status text to "In Translation" on every
# sheet that references it (Overview: E2, F2 ; zh-cn: C2 ; de-de: C2).
foreach ($ws in $wb.Worksheets) {
    $used = $ws.UsedRange
    foreach ($cell in $used.Cells) {
        $text = $cell.Text.ToString()
        if ($text -eq "Ready for handoff") {
            $cell.Value = "In Translation"
        }
    }
}

# The status column is narrower now that the text is shorter than before
# ("Ready for handoff" -> "In Translation"); re-fit those columns so the
# stored column width reflects the new, shorter text, same as Excel would
# do automatically when the cell content changes.
$newWidth = 12.5

$overview = $wb.Worksheets.Item("Overview")
$overview.Columns.Item(5).ColumnWidth = $newWidth
$overview.Columns.Item(6).ColumnWidth = $newWidth

$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Columns.Item(3).ColumnWidth = $newWidth

$dede = $wb.Worksheets.Item("de-de")
$dede.Columns.Item(3).ColumnWidth = $newWidth
